$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2.94
$ws.Range("G2").Value = 3.25
$ws.Range("H2").Value = 2.68
$ws.Range("I2").Value = 2.94
$ws.Range("J2").Value = 3.1
